# Commit: "updated demo description with --dry-run"
#
# The edited content lives on the *Notes Page* of slide 13 ("Demo"),
# inside the notes body placeholder ("Notes Placeholder 2"). Two
# changes are made there:
#   1. A new bullet "When creating pods, demo the --dry-run flag" is
#      inserted right before the "Create a pod: ~/..." bullet.
#   2. The demo file referenced in that next bullet is renamed from
#      02_pod_exec_liveness.yaml to 02b_pod_exec_liveness.yaml.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(13)
$notesPage = $slide.NotesPage

# Locate the notes body placeholder robustly by name rather than a
# hard-coded index.
$notesShape = $null
for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
    $candidate = $notesPage.Shapes.Item($i)
    if ($candidate.Name -eq "Notes Placeholder 2") {
        $notesShape = $candidate
    }
}

$notesShape.TextFrame.TextRange.Text = "Start with the “kubectl explain pod” and “kubectl explain pod.spec”`nShow how to get an overview as well as detailed info about a resource type.`nWhen creating pods, demo the --dry-run flag`nCreate a pod: ~/kubernetes/demo/02b_pod_exec_liveness.yaml`nDiscuss the probe and how it should fail`nShow how it fails & get restarted`nPoint out the failure threshold`nCreate a 2nd pod, this time with a web server: ~/kubernetes/demo/02_pod_http_liveness.yaml`nExplain the http probe and how it should fill up the logs`nShow logs of the container (will be the access log) and discuss the effect of the liveness probe`nAccess nginx:`nRun  kubectl port-forward pod/nginx-liveness-pod 8080:80`nOpen a browser and connect to 127.0.0.1:8080`nPort-forward is a nice command to test access to something that you don’t want to expose (yet). However it is not recommended for any production like setup as the traffic is routed via the cluster’s API server"
